# Update "想去人数" (F column) and one "最低票价" (G column) values on the
# "展览" and "全部类型" worksheets, as produced by the latest data refresh.

$wb = $excel.ActiveWorkbook

# Row number -> @(Column, NewValue)
$updates = @{
    2  = @("F", 25)
    4  = @("F", 119)
    5  = @("F", 22)
    7  = @("F", 88)
    8  = @("F", 461)
    11 = @("F", 573)
    13 = @("F", 300)
    15 = @("F", 372)
    17 = @("F", 91)
    19 = @("F", 51)
    21 = @("F", 97)
    22 = @("F", 935)
    23 = @("F", 1397)
    24 = @("F", 300)
    26 = @("F", 184)
    27 = @("F", 75)
    28 = @("F", 154)
    29 = @("G", 55)
    32 = @("F", 250)
    33 = @("F", 275)
    34 = @("F", 1621)
    35 = @("F", 52)
    38 = @("F", 583)
    40 = @("F", 3650)
    42 = @("F", 202)
    43 = @("F", 911)
    44 = @("F", 42)
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $col = $updates[$row][0]
        $val = $updates[$row][1]
        $ws.Range("$col$row").Value = $val
    }
}

$wb.Save()
